# Team-Meeting-Info.xlsx update:
# Add four new meeting-log entries (dates 44165, 44166, 44168, 44170) to the
# meeting table on Sheet1, widen column C to fit the new (longer) text,
# and update the sheet's active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Copy-Format($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial(-4122) | Out-Null
}

# ---------------------------------------------------------------------
# First lay down the formatting (borders / alignment) for every new row,
# mirroring the existing table rows above them.
# ---------------------------------------------------------------------
Copy-Format "A34:D34" "A63:D63"

Copy-Format "A58" "A64"
Copy-Format "B58" "B64"
Copy-Format "C58" "C64"
Copy-Format "D58" "D64"

Copy-Format "A34:D34" "A65:D66"
$ws.Range("A65:A66").HorizontalAlignment = -4131

Copy-Format "A58" "A67"
Copy-Format "B58" "B67"
Copy-Format "C58" "C67"
Copy-Format "D58" "D67"

Copy-Format "A34:D34" "A68:D69"
$ws.Range("A68:A69").HorizontalAlignment = -4131

Copy-Format "A58" "A70"
Copy-Format "B58" "B70"
Copy-Format "C58" "C70"
Copy-Format "D58" "D70"

Copy-Format "A34:D34" "A71:D71"
$ws.Range("A71").HorizontalAlignment = -4131

Copy-Format "A34:D34" "A72:D73"
$ws.Range("A72:A73").HorizontalAlignment = -4131

Copy-Format "A58" "A74"
Copy-Format "B58" "B74"
Copy-Format "C58" "C74"
Copy-Format "D58" "D74"

Copy-Format "A34:D34" "A75:D75"
Copy-Format "A34:D34" "A76:D76"

# ---------------------------------------------------------------------
# Now fill in the column-C "Meeting Highlights" text, top to bottom.
# ---------------------------------------------------------------------
$ws.Range("C64").Value = "The team gathered to discuss more on code analysis, tools, and task segregation."
$ws.Range("C67").Value = "Discussion on manual and automated code analysis."
$ws.Range("C70").Value = "Team gathered to discuss more on the feedback received from team check-in with the professor"
$ws.Range("C71").Value = "Set proper code review strategy and related CWE’s are consolidated for further analysis."
$ws.Range("C74").Value = "Team discussed on their respective manual and automated code analysis progress."
$ws.Range("C75").Value = "Peer review and refined some of the code review strategies for further analysis."
$ws.Range("C76").Value = "Other than all these meetings team was actively available on WhatsApp as it was a big assignment and discussed any blockage on the assigned tasks then and there."

# ---------------------------------------------------------------------
# Fill in the column-B "Time of meeting" text, top to bottom.
# ---------------------------------------------------------------------
$ws.Range("B64").Value = " 7:30 to 9 PM"
$ws.Range("B67").Value = " 7:30 to 9 PM"
$ws.Range("B70").Value = "7:00 to 8:30 PM"
$ws.Range("B74").Value = "12:00 to 1:00PM"

# ---------------------------------------------------------------------
# Fill in column-A meeting dates.
# ---------------------------------------------------------------------
$ws.Range("A64").Value = 44165
$ws.Range("A67").Value = 44166
$ws.Range("A70").Value = 44168
$ws.Range("A74").Value = 44170

# ---------------------------------------------------------------------
# Fill in column-D attendance.
# ---------------------------------------------------------------------
$ws.Range("D64").Value = "Full"
$ws.Range("D67").Value = "Full"
$ws.Range("D70").Value = "Full"
$ws.Range("D74").Value = "Full"

# ---------------------------------------------------------------------
# Column C needs to be wide enough to show the longer text added above.
# ---------------------------------------------------------------------
$ws.Columns("C").ColumnWidth = 138.42857142857142

# ---------------------------------------------------------------------
# Update the sheet selection to match the edited range.
# ---------------------------------------------------------------------
$ws.Range("A10:D76").Select() | Out-Null

Write-Output "done"
